$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look numeric,
# so Excel keeps them as text (preserving formatting/trailing zeros)
# instead of silently converting to a floating point number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "72.238.60"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "4.036.82"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "540.08"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "151.81"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "4.029.81"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "0.753"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").Value = "0.172"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "53.84"
$ws.Range("E12").Value = "  +11.16%  "
$ws.Range("D13").Value = "0.0000330"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "4.682.39"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "4.038.98"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "14.34"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "20.68"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "72.169.23"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "450.80"
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("D24").Value = "3.52"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").Value = "4.26"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "14.59"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "4.29"
$ws.Range("E27").Value = "  +15.51%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").Value = "5.95"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "8.02"
$ws.Range("E32").Value = "  +16.57%  "
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").Value = "13.60"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").Value = "49.37"
$ws.Range("E35").Value = "  +15.82%  "
$ws.Range("D36").Value = "682.46"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "66.75"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "0.450"
$ws.Range("E38").Value = "  +4.83%  "
$ws.Range("D39").Value = "0.0₃0889"
$ws.Range("E39").Value = "  +5.48%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "3.47"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.148"
$ws.Range("E41").Value = "  -5.71%  "
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").Value = "11.18"
$ws.Range("E43").Value = "  +16.63%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "0.0493"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "2.66"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").Value = "3.11"
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("D50").Value = "3.32"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "3.47"
$ws.Range("E51").Value = "  +3.35%  "

# Restore default style on the forced-text cells
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
